# Report generation: when a sale-phu report has no data rows, the export
# script still writes one blank/zero "template" row right under the header
# row so the downstream formatting (borders, totals, number formats) has
# somewhere to land. Replicate that for the "Đơn sale phụ" sheet: add row 2
# with empty text columns and 0 in every money column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn sale phụ")

# Text / identifier columns on row 2 stay blank.
$textCols = @("A", "C", "D", "E", "F", "G", "H", "J", "Q", "R", "S", "T")
foreach ($col in $textCols) {
    $ws.Range($col + "2").Value = ""
}

# Mã dịch vụ (B) is a numeric field with no value for the blank row.
$ws.Range("B2").Value = ""

# Money / numeric columns on row 2 default to 0.
$numCols = @("I", "K", "L", "M", "N", "O", "P")
foreach ($col in $numCols) {
    $ws.Range($col + "2").Value = 0
}
